# Restored from revision of admin on 04/14/2021 12:02:04 PM.TEST
# Author: admin. Type: SAVE.
#
# The only meaningful content change in this revision is the "Integer min"
# rule value on the Rules sheet (cell C10), which changes from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
